$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Reorganizacion de las tablas grupo-clases ---

# Row 16 (student #6) currently points at group "03" class entries
# (Practicas de Aula/Semin-02 / Practicas de Laboratorio-03 / Tutorias
# Grupales-03). Re-point it to the "01" group, same as row 12.
$ws.Range("J16").Value = "Practicas de Aula/Semin-01"
$ws.Range("K16").Value = "Prácticas de Laboratorio-01"
$ws.Range("L16").Value = "Tutorías Grupales-01"

# Row 17 held student #8 (Gallego Doncel, Aljenadro) together with a
# mailto hyperlink on the email cell. Remove the hyperlink first, then
# clear the row contents while keeping the existing cell formatting.
$ws.Range("D17").Hyperlinks.Delete()
$ws.Range("A17:L17").ClearContents()

# Update the active selection left in the sheet so it now highlights the
# freshly emptied row instead of the old I17:L17 block.
$ws.Range("A17:XFD17").Select()

$wb.Save()
